$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.198.74"
$ws.Range("D3").Value = "2.484.19"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.78"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.06"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.64%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.15"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +4.00%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.38"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.17"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "2.871.91"
$ws.Range("E15").Value = "  -0.48%  "
$ws.Range("D16").Value = "2.483.98"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.843"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "47.102.71"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.41"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.61"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "0.0₃0943"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.74"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +14.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.51"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "245.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -2.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.55"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +3.65%  "
$ws.Range("E29").Value = "  -1.24%  "
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.98%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.76"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.24"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.21%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.33"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0781"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.74"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.93"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.05"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +8.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.112"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.26%  "
$ws.Range("E42").Value = "  -0.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "118.07"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.65%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").Value = "1.995.44"
$ws.Range("E45").Value = "  +1.31%  "
$ws.Range("E46").Value = "  +1.37%  "
$ws.Range("E47").Value = "  -3.89%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.79"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("D49").ClearFormats()
$ws.Range("E50").Value = "  -4.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "56.48"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.84%  "
